# Generate Report for Handback
# Adds a new handback row (5a454809-12c8-4e4b-a183-ea6968e68ac6) to the
# Overview sheet and to each per-locale sheet (zh-cn, de-de), mirroring the
# existing rows 2/3 pattern, then grows the tables/dimensions accordingly.

$wb = $excel.ActiveWorkbook

$dateFmt = "yyyy-mm-dd HH:mm:ss"
$hyperlinkColor = 15570276  # OLE BGR for RGB FF6495ED (matches existing HyperLink font)

# ---------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A4").Value = "5a454809-12c8-4e4b-a183-ea6968e68ac6.md"
$wsOverview.Range("C4").Value = ".md"
$wsOverview.Range("E4").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F4").Value = "Handed back: in sync with en-US"

$wsOverview.Range("G4").NumberFormat = $dateFmt
$wsOverview.Range("G4").Value = "'2016-08-19 18:53:10"

$wsOverview.Range("B4").Font.Underline = $true
$wsOverview.Range("B4").Font.Color = $hyperlinkColor
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bb3a204ecbb8096af1ff356de916c08ddd2cc84a/e2e/5a454809-12c8-4e4b-a183-ea6968e68ac6.md", "", "", "e2e\5a454809-12c8-4e4b-a183-ea6968e68ac6.md")

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ---------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("B4").Value = ".md"
$wsZhCn.Range("C4").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("D4").Value = "e2e"
$wsZhCn.Range("E4").Value = "ht"
$wsZhCn.Range("F4").Value = "'True"

$wsZhCn.Range("G4").Value = "5a454809-12c8-4e4b-a183-ea6968e68ac6.bb3a204ecbb8096af1ff356de916c08ddd2cc84a.zh-cn.xlf"
$wsZhCn.Range("H4").NumberFormat = $dateFmt
$wsZhCn.Range("H4").Value = "'2016-08-19 18:53:02"

$wsZhCn.Range("J4").Value = "5a454809-12c8-4e4b-a183-ea6968e68ac6.bb3a204ecbb8096af1ff356de916c08ddd2cc84a.zh-cn.xlf"
$wsZhCn.Range("K4").NumberFormat = $dateFmt
$wsZhCn.Range("K4").Value = "'2016-08-19 18:53:28"

$wsZhCn.Range("L4").Value = "'"
$wsZhCn.Range("M4").Value = "'True"
$wsZhCn.Range("N4").Value = "'"
$wsZhCn.Range("O4").Value = "'False"
$wsZhCn.Range("P4").Value = "'"

$wsZhCn.Range("A4").Font.Underline = $true
$wsZhCn.Range("A4").Font.Color = $hyperlinkColor
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bb3a204ecbb8096af1ff356de916c08ddd2cc84a/e2e/5a454809-12c8-4e4b-a183-ea6968e68ac6.md", "", "", "5a454809-12c8-4e4b-a183-ea6968e68ac6.md")

$wsZhCn.Range("I4").Font.Underline = $true
$wsZhCn.Range("I4").Font.Color = $hyperlinkColor
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/bb3a204ecbb8096af1ff356de916c08ddd2cc84a/e2e/5a454809-12c8-4e4b-a183-ea6968e68ac6.md", "", "", "5a454809-12c8-4e4b-a183-ea6968e68ac6.md")

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ---------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("B4").Value = ".md"
$wsDeDe.Range("C4").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("D4").Value = "e2e"
$wsDeDe.Range("E4").Value = "ht"
$wsDeDe.Range("F4").Value = "'True"

$wsDeDe.Range("G4").Value = "5a454809-12c8-4e4b-a183-ea6968e68ac6.bb3a204ecbb8096af1ff356de916c08ddd2cc84a.de-de.xlf"
$wsDeDe.Range("H4").NumberFormat = $dateFmt
$wsDeDe.Range("H4").Value = "'2016-08-19 18:53:10"

$wsDeDe.Range("J4").Value = "5a454809-12c8-4e4b-a183-ea6968e68ac6.bb3a204ecbb8096af1ff356de916c08ddd2cc84a.de-de.xlf"
$wsDeDe.Range("K4").NumberFormat = $dateFmt
$wsDeDe.Range("K4").Value = "'2016-08-19 18:53:35"

$wsDeDe.Range("L4").Value = "'"
$wsDeDe.Range("M4").Value = "'True"
$wsDeDe.Range("N4").Value = "'"
$wsDeDe.Range("O4").Value = "'False"
$wsDeDe.Range("P4").Value = "'"

$wsDeDe.Range("A4").Font.Underline = $true
$wsDeDe.Range("A4").Font.Color = $hyperlinkColor
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bb3a204ecbb8096af1ff356de916c08ddd2cc84a/e2e/5a454809-12c8-4e4b-a183-ea6968e68ac6.md", "", "", "5a454809-12c8-4e4b-a183-ea6968e68ac6.md")

$wsDeDe.Range("I4").Font.Underline = $true
$wsDeDe.Range("I4").Font.Color = $hyperlinkColor
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/bb3a204ecbb8096af1ff356de916c08ddd2cc84a/e2e/5a454809-12c8-4e4b-a183-ea6968e68ac6.md", "", "", "5a454809-12c8-4e4b-a183-ea6968e68ac6.md")

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Output "Handback row added for 5a454809-12c8-4e4b-a183-ea6968e68ac6"
